$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the D1/E1 header text to reflect the new combined LS1088a+LS2088a naming
$ws.Range("D1").Value = "LS1088a Rev1 + LS2088a Rev1 Implementation"
$ws.Range("E1").Value = "LS1088a Rev1 + LS2088a Rev1 API"

# Resize columns D/E to fit the new (longer) header text
$ws.Columns.Item(4).ColumnWidth = 40.65
$ws.Columns.Item(5).ColumnWidth = 29.3

# Move the selection to D1, matching the author's final cursor position
$ws.Activate()
$ws.Range("D1").Select() | Out-Null
